# feat: Se agrego paginacion
# Refresh the existing listing rows (2-11) with a new scrape and append a second
# results page (rows 12-31) using the same title/price/list_currency/rating/site columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding numeric-looking text (e.g. "21,999", "5.0") need a leading quote so
# Excel keeps them as literal text instead of silently converting them to numbers.
function Set-CellText($cell, [string]$text) {
    if ($text -match '^[0-9][0-9,\.]*$') {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

$rows = @(
    @("HP Laptop 14-ew1000la Intel Core Ultra 5 16GB RAM 512GB SSD", "21,999", "`$", "N/A", "mercadolibre"),
    @("Notebook Asus Vivobook Intel Core I3 1215u, 8gb De Ram, 128gb Ssd, Full Hd, Windows 11 Home", "6,190", "`$", "N/A", "mercadolibre"),
    @("Laptop Acer Aspire 3 15.6 Ryzen 7, 16gb/512gb, Windows 11 Color Plateado", "7,999", "`$", "4.5", "mercadolibre"),
    @("Laptop 15.6in Notebook Intel Celeron N5095 16g Ram 256g Ssd Color Plateado", "5,299", "`$", "N/A", "mercadolibre"),
    @("Notebook N5095 15.6`" gris 8GB de Ram - 256GB SSD - Intel Celeron", "6,941", "`$", "N/A", "mercadolibre"),
    @("Laptop Lenovo Ideapad Celeron 4gb + 128ssd + Office Regalo Color Gris", "3,699", "`$", "4.6", "mercadolibre"),
    @("Apple MacBook Air (13 pulgadas, 2020, Chip M1, 256 GB de SSD, 8 GB de RAM) - Gris espacial - Distribuidor Autorizado", "11,999", "`$", "4.9", "mercadolibre"),
    @("Laptop gamer Thunderobot 911MT 12th Intel Core i7 12650H 16GB de RAM 512GB SSD, NVIDIA GeForce RTX 3050 165 Hz 1920x1080px Windows 11 Pro", "15,762", "`$", "4.7", "mercadolibre"),
    @("Laptop Acer Etbook Intel Core I5-12450h 16gb 512gb Ssd Win11", "8,618", "`$", "N/A", "mercadolibre"),
    @("Portátil Vanwin 11.6 Intel 8gb Ram 256g Rom Windows 11 Gris", "3,855", "`$", "4.2", "mercadolibre"),
    @("HP Laptop 14-ew1000la Intel Core Ultra 5 16GB RAM 512GB SSD", "21,999", "`$", "N/A", "mercadolibre"),
    @("Laptop Acer Etbook Intel Core I5-12450h 16gb 512gb Ssd Win11", "8,618", "`$", "N/A", "mercadolibre"),
    @("Portátil Vanwin 11.6 Intel 8gb Ram 256g Rom Windows 11 Gris", "3,855", "`$", "4.2", "mercadolibre"),
    @("Laptop Lenovo Ideapad Slim 3 15.6'' Ci5 8gb + 512gb Ssd", "8,799", "`$", "4.8", "mercadolibre"),
    @("Ultrabook Win VNJH1601-3 16`" plateado 16GB de Ram - 1TB SSD - Intel Core i7", "13,635", "`$", "N/A", "mercadolibre"),
    @("Lenovo T490 Thinkpad Core I5 8va 8gb 256 Ssd 14 Fhd W10 Pro", "4,317", "`$", "4.6", "mercadolibre"),
    @("Notebook V-Series V14 14`" black 4GB de Ram - 128GB SSD - Intel Celeron", "3,801", "`$", "4.6", "mercadolibre"),
    @("Apple MacBook Air (13 pulgadas, 2020, Chip M1, 256 GB de SSD, 8 GB de RAM) - Gris espacial", "12,199", "`$", "4.9", "mercadolibre"),
    @("Laptop Moderness Lt1505 16+512g 15.6 1920*1080p Windows 11", "5,499", "`$", "4.9", "mercadolibre"),
    @("Notebook Asus Vivobook Intel Core I3 1215u, 8gb De Ram, 128gb Ssd, Full Hd, Windows 11 Home", "6,190", "`$", "N/A", "mercadolibre"),
    @("Ordenador Portátil Windows 11 Intel Core I7 16gb+1tb Ssd", "10,566", "`$", "N/A", "mercadolibre"),
    @("Notebook N5095 15.6`" gris 8GB de Ram - 256GB SSD - Intel Celeron", "6,941", "`$", "N/A", "mercadolibre"),
    @("Laptop 15.6inch Intel 2 En 1 Ddr4 De 16gb Ram 512 Gb", "5,308", "`$", "N/A", "mercadolibre"),
    @("Notebook Stream 14-cf2112wm rosa 4GB de Ram - 64GB SSD - Intel Celeron", "4,211", "`$", "4.4", "mercadolibre"),
    @("Notebook Acer Aspire 3 Ryzen 5 7520u 512g 16g 15.6 Touch Color Gris", "10,036", "`$", "5.0", "mercadolibre"),
    @("Laptop Lenovo Ideapad Slim 3 Core I7 Ram 16gb Ssd 1tb W11h Color Gris", "16,799", "`$", "4.7", "mercadolibre"),
    @("Laptop Lenovo Ideapad 15.6 Ryzen 3 7320u 8gb 256gb ssd", "6,799", "`$", "4.8", "mercadolibre"),
    @("Laptop Asus 15-x5ja Core I3-1005, 512gb Ssd 8gb Ram, Fhd", "5,871", "`$", "N/A", "mercadolibre"),
    @("Notebook Msi Amd Ryzen 5-7530u 16gb 512gb 14'' Fhd Black Color Negro", "11,399", "`$", "4.6", "mercadolibre"),
    @("Yinkuu Laptop Intel Celeron 16gb+512gb Ssd 15.6 Uhd 4k", "5,129", "`$", "N/A", "mercadolibre"),
)

$r = 2
foreach ($rowVals in $rows) {
    for ($col = 1; $col -le 5; $col++) {
        Set-CellText $ws.Cells.Item($r, $col) $rowVals[$col - 1]
    }
    $r++
}
